# Transfermarkt World Cup 2022 Comparison.xlsx
# "Rerun all experiment : convert linear regression to logistic regression"
#
# This script reproduces the re-run of the ranking experiment: the
# "R-Proposed" values (column D) on the three "Top 10 players ..." sheets
# change (new model scores), a few Win/Lose/Draw results flip, and the
# "Top 10 proposed ranking" sheet is reshuffled with two players replaced
# by new ones. The previously-active tab ("Top 10 proposed ranking") is
# swapped for "Top 10 VAEP ranking".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: "Top 10 players goal 90"  (R-Proposed column D, Result column E)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Top 10 players goal 90")
$ws.Range("D2").Value = 125
$ws.Range("D3").Value = 150
$ws.Range("D4").Value = 107
$ws.Range("D5").Value = 102
$ws.Range("D6").Value = 153
$ws.Range("D7").Value = 84
$ws.Range("D8").Value = 100
$ws.Range("E8").Value = "Lose"
$ws.Range("D9").Value = 165
$ws.Range("D10").Value = 135
$ws.Range("D11").Value = 158
$ws.Range("D2:D11").Select()

# ---------------------------------------------------------------
# Sheet: "Top 10 players assist 90"  (R-Proposed column D, Result column E)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Top 10 players assist 90")
$ws.Range("D2").Value = 84
$ws.Range("D4").Value = 35
$ws.Range("E4").Value = "Draw"
$ws.Range("D5").Value = 150
$ws.Range("D6").Value = 66
$ws.Range("D7").Value = 88
$ws.Range("D8").Value = 67
$ws.Range("D9").Value = 90
$ws.Range("D10").Value = 100
$ws.Range("E10").Value = "Lose"
$ws.Range("D11").Value = 102
$ws.Range("D2:D11").Select()

# ---------------------------------------------------------------
# Sheet: "Top 10 players goal assist 90"  (R-Proposed column D, Result column E)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Top 10 players goal assist 90")
$ws.Range("D2").Value = 84
$ws.Range("D3").Value = 150
$ws.Range("D4").Value = 125
$ws.Range("D5").Value = 107
$ws.Range("D6").Value = 35
$ws.Range("E6").Value = "Draw"
$ws.Range("D8").Value = 102
$ws.Range("D9").Value = 153
$ws.Range("D10").Value = 100
$ws.Range("E10").Value = "Lose"
$ws.Range("D11").Value = 66

# ---------------------------------------------------------------
# Sheet: "Top 10 proposed ranking"  (players reshuffled, two replaced)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Top 10 proposed ranking")
$ws.Range("B5").Value = "Carlos Soler"
$ws.Range("C5").Value = "50 million euro"
$ws.Range("B6").Value = "Jorge Resurreccion Merodio "
$ws.Range("C6").Value = "35 million euro"
$ws.Range("B8").Value = "Ilkay Gundogan"
$ws.Range("C8").Value = "25 million euro"
$ws.Range("B9").Value = "Niklas Sule"
$ws.Range("C9").Value = "35 million euro"
$ws.Range("B10").Value = "Toby Aldeweireld"
$ws.Range("C10").Value = "8 million euro"
$ws.Range("B11").Value = "Hattan Bahebri"
$ws.Range("C11").Value = "0,7 million euro"
$ws.Range("C13").Value = "253,5 million euro"
$ws.Range("G16").Select()

# ---------------------------------------------------------------
# Make "Top 10 VAEP ranking" the active tab (was "Top 10 proposed ranking")
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Top 10 VAEP ranking")
$ws.Activate()
